# Update TPM-derived NATMI ligand-receptor metrics (Col18a1-Itga5)
# Updates columns E..T for rows 2..10 with new values as computed with updated TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5586043333333334
$ws.Range("H2").Value = 1.675813
$ws.Range("I2").Value = 0.01643366487114074
$ws.Range("J2").Value = 0.01643366487114074
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 5.780221833859222
$ws.Range("R2").Value = 52.021996504733
$ws.Range("S2").Value = 0.003779666467045227
$ws.Range("T2").Value = 0.003779666467045227
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5586043333333334
$ws.Range("H3").Value = 1.675813
$ws.Range("I3").Value = 0.01643366487114074
$ws.Range("J3").Value = 0.01643366487114074
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 17.07645756880344
$ws.Range("R3").Value = 153.688118119231
$ws.Range("S3").Value = 0.01116623477504739
$ws.Range("T3").Value = 0.01116623477504739
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5586043333333334
$ws.Range("H4").Value = 1.675813
$ws.Range("I4").Value = 0.01643366487114074
$ws.Range("J4").Value = 0.01643366487114074
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 2.275228221120889
$ws.Range("R4").Value = 20.477053990088
$ws.Range("S4").Value = 0.001487763629048122
$ws.Range("T4").Value = 0.001487763629048122
$ws.Range("H5").Value = 63.825936
$ws.Range("I5").Value = 0.6259016025719319
$ws.Range("J5").Value = 0.6259016025719319
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 220.1487092137973
$ws.Range("R5").Value = 1981.338382924176
$ws.Range("S5").Value = 0.1439544567484407
$ws.Range("T5").Value = 0.1439544567484408
$ws.Range("H6").Value = 63.825936
$ws.Range("I6").Value = 0.6259016025719319
$ws.Range("J6").Value = 0.6259016025719319
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("Q6").Value = 650.3833589387147
$ws.Range("S6").Value = 0.4252833616359039
$ws.Range("T6").Value = 0.425283361635904
$ws.Range("H7").Value = 63.825936
$ws.Range("I7").Value = 0.6259016025719319
$ws.Range("J7").Value = 0.6259016025719319
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 86.65559392763733
$ws.Range("R7").Value = 779.900345348736
$ws.Range("S7").Value = 0.05666378418758725
$ws.Range("T7").Value = 0.05666378418758726
$ws.Range("G8").Value = 12.157548
$ws.Range("H8").Value = 36.472644
$ws.Range("I8").Value = 0.3576647325569273
$ws.Range("J8").Value = 0.3576647325569273
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 125.801609837956
$ws.Range("R8").Value = 1132.214488541604
$ws.Range("S8").Value = 0.08226122454669958
$ws.Range("T8").Value = 0.08226122454669961
$ws.Range("G9").Value = 12.157548
$ws.Range("H9").Value = 36.472644
$ws.Range("I9").Value = 0.3576647325569273
$ws.Range("J9").Value = 0.3576647325569273
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 371.654568670892
$ws.Range("R9").Value = 3344.891118038028
$ws.Range("S9").Value = 0.2430235985582661
$ws.Range("T9").Value = 0.2430235985582661
$ws.Range("G10").Value = 12.157548
$ws.Range("H10").Value = 36.472644
$ws.Range("I10").Value = 0.3576647325569273
$ws.Range("J10").Value = 0.3576647325569273
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 49.518406246816
$ws.Range("R10").Value = 445.665656221344
$ws.Range("S10").Value = 0.03237990945196165
$ws.Range("T10").Value = 0.03237990945196165
